$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.195.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.909.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.10%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5085'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.48%  '

$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09338'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.139'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.48%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.91'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.404'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.89'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.908.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.324'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001122'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06614'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("E20").Value = '  +1.82%  '

$ws.Range("E21").Value = '  -0.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.231'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.251.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.50'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.324'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.596'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.121.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.60%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '157.96'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.105'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1075'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.650'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.40%  '

$ws.Range("E34").Value = '  -0.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.723'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.33%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06688'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02426'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.64%  '

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2208'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.78%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.245'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.23%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.285'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.63%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6534'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.85%  '

$ws.Range("E42").Value = '  +0.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.018'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6122'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.32'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.726'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.82%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.288'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.023'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.94%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '123.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.83%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.188'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.70%  '
